# Update Name of Algo - apply revised values to result_data_RandomForest.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.40969999999999
$ws.Range("E6").Value = 12.16519999999999
$ws.Range("E7").Value = 12.10859999999999
$ws.Range("C8").Value = -12.2724
$ws.Range("E8").Value = 13.6685
$ws.Range("A12").Value = -21.88090000000002
$ws.Range("C12").Value = -12.7008
$ws.Range("C14").Value = -11.99859999999999
$ws.Range("E19").Value = 12.9715
$ws.Range("E21").Value = 12.67299999999999
$ws.Range("C22").Value = -11.07189999999999
$ws.Range("E24").Value = 12.64239999999999
